$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("展览").Range("F3").Value = 168
$wb.Worksheets.Item("展览").Range("F5").Value = 7
$wb.Worksheets.Item("展览").Range("F6").Value = 338
$wb.Worksheets.Item("展览").Range("F7").Value = 5948
$wb.Worksheets.Item("展览").Range("F9").Value = 4006
$wb.Worksheets.Item("展览").Range("F13").Value = 128
$wb.Worksheets.Item("展览").Range("F15").Value = 3993
$wb.Worksheets.Item("展览").Range("F16").Value = 13
$wb.Worksheets.Item("展览").Range("F19").Value = 5587
$wb.Worksheets.Item("展览").Range("F21").Value = 2200
$wb.Worksheets.Item("展览").Range("F22").Value = 147
$wb.Worksheets.Item("展览").Range("F23").Value = 395
$wb.Worksheets.Item("展览").Range("F24").Value = 8373
$wb.Worksheets.Item("展览").Range("F26").Value = 42
$wb.Worksheets.Item("展览").Range("F27").Value = 2239
$wb.Worksheets.Item("展览").Range("F28").Value = 2276
$wb.Worksheets.Item("展览").Range("F29").Value = 1353
$wb.Worksheets.Item("展览").Range("F31").Value = 1859
$wb.Worksheets.Item("展览").Range("F32").Value = 36
$wb.Worksheets.Item("展览").Range("F33").Value = 295
$wb.Worksheets.Item("展览").Range("F44").Value = 1387
$wb.Worksheets.Item("展览").Range("F45").Value = 2234
$wb.Worksheets.Item("演出").Range("F3").Value = 2
$wb.Worksheets.Item("本地生活").Range("F2").Value = 628
$wb.Worksheets.Item("本地生活").Range("F3").Value = 822
$wb.Worksheets.Item("全部类型").Range("F3").Value = 168
$wb.Worksheets.Item("全部类型").Range("F4").Value = 628
$wb.Worksheets.Item("全部类型").Range("F5").Value = 822
$wb.Worksheets.Item("全部类型").Range("F6").Value = 338
$wb.Worksheets.Item("全部类型").Range("F7").Value = 5948
$wb.Worksheets.Item("全部类型").Range("F9").Value = 4006
$wb.Worksheets.Item("全部类型").Range("F13").Value = 128
$wb.Worksheets.Item("全部类型").Range("F17").Value = 3993
$wb.Worksheets.Item("全部类型").Range("F18").Value = 13
$wb.Worksheets.Item("全部类型").Range("F21").Value = 5587
$wb.Worksheets.Item("全部类型").Range("F23").Value = 2200
$wb.Worksheets.Item("全部类型").Range("F24").Value = 147
$wb.Worksheets.Item("全部类型").Range("F25").Value = 395
$wb.Worksheets.Item("全部类型").Range("F26").Value = 8373
$wb.Worksheets.Item("全部类型").Range("F29").Value = 2239
$wb.Worksheets.Item("全部类型").Range("F30").Value = 2276
$wb.Worksheets.Item("全部类型").Range("F31").Value = 1353
$wb.Worksheets.Item("全部类型").Range("F33").Value = 1859
$wb.Worksheets.Item("全部类型").Range("F34").Value = 36
$wb.Worksheets.Item("全部类型").Range("F35").Value = 295
$wb.Worksheets.Item("全部类型").Range("F45").Value = 1387
$wb.Worksheets.Item("全部类型").Range("F46").Value = 2234
